# Master data refresh (16th May) — append 3 new user_detail_h rows
# and fix a pre-existing style inconsistency on I33.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix I33: it was missing the "left aligned" style that every other
#     row in the is_active column carries (I2:I32, I34:I36) ---
$ws.Range("I33").HorizontalAlignment = -4131

# --- New rows 34-36 ---
$newRows = @(
    @{ Row = 34; Id = 110033; Uin = 9317596771; Name = "Nikola Tesla"; Email = "nikola.tesla@xyz.com"; Mobile = 818876434 },
    @{ Row = 35; Id = 110034; Uin = 9317596772; Name = "Graham Bell";  Email = "graham.bell@xyz.com";  Mobile = 818876435 },
    @{ Row = 36; Id = 110035; Uin = 9317596773; Name = "Albert Miles"; Email = "albert.miles@xyz.com"; Mobile = 818876436 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $r.Id
    $ws.Cells.Item($row, 2).Value = $r.Uin
    $ws.Cells.Item($row, 3).Value = $r.Name
    $ws.Cells.Item($row, 4).Value = $r.Email
    $ws.Cells.Item($row, 5).Value = $r.Mobile
    $ws.Cells.Item($row, 6).Value = "ACT"
    $ws.Cells.Item($row, 7).Value = "eng"
    $ws.Cells.Item($row, 8).Value = "PWD"
    $ws.Cells.Item($row, 9).Value = $true
    $ws.Cells.Item($row, 10).Value = "superadmin"
    $ws.Cells.Item($row, 11).Value = "now()"
    $ws.Cells.Item($row, 12).Value = "now()"

    # Column I (is_active) carries the same left-aligned style as the rest
    # of the column; column D (email) carries the fill-applied style that
    # the rest of the column (D2:D32) carries.
    $ws.Cells.Item($row, 9).HorizontalAlignment = -4131
    $ws.Cells.Item($row, 4).Style = $ws.Cells.Item(32, 4).Style
}

# --- Selection cursor moved from M6 back to M1 (still the whole M:XFD block) ---
$null = $ws.Range("M1:XFD1048576").Select()
